$d = $word.ActiveDocument

# Update model prediction / coefficient analysis table values.
$d.Content.Find.Execute("8.874", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "39.470", 2)

$d.Content.Find.Execute("0.012", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "0.000", 2)

$d.Content.Find.Execute("1.966", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2.843", 2)

$d.Content.Find.Execute("0.161", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "0.092", 2)

$d.Content.Find.Execute("0.463", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "11.642", 2)

$d.Content.Find.Execute("0.793", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "0.003", 2)
